# Fixed capabilities sync for future and dynamic items, fixed tier1 and
# tier2 templates sync.
#
# The "Templates" sheet gets two new placeholder rows inserted right after
# the existing header/first-template row: one default activation template
# row for scope "tier1" and one for scope "tier2". The previously-existing
# template rows shift down by two. The "Templates" sheet also becomes the
# active sheet/selection (it had been left on "Configuration").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Templates")

# Insert two new rows above the old row 3, pushing the existing rows 3-5
# down to rows 5-7. Insert() duplicates the formatting of the row above,
# which also grows the data validations (C/D/E columns) automatically.
$ws.Rows("3:4").Insert()

# New row 3: default activation template for tier1
$ws.Range("A3").Value = "TL-559-508-354"
$ws.Range("B3").Value = "Default Activation Template"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "tier1"
$ws.Range("E3").Value = "fulfillment"

# New row 4: default activation template for tier2
$ws.Range("A4").Value = "TL-518-222-757"
$ws.Range("B4").Value = "Default Activation Template"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "tier2"
$ws.Range("E4").Value = "fulfillment"

# Move the active tab/selection to Templates (it had been parked on the
# last sheet, "Configuration").
$ws.Activate()
$ws.Range("E3:E4").Select()
